# PSX balance sheet: the "current period" column (B) had been left as a
# column of empty placeholders (inline-string "") while every other period
# column (C, D, E, ...) already carried real figures. This fills column B
# in with its reported numbers and widens it to line up with the other
# (already 16.5-wide) data columns instead of its old, narrower 12.1 width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PSX")

# Column B was narrower than the other data columns (bestFit from when it
# held only blanks). Match it to column C's width now that it holds real
# numbers too.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Row -> column B value, keyed by the row's label in column A.
$rowValues = @{
    2  = 1351000000.0    # Cash and Short Term Investments
    3  = 8165000000.0    # Receivables
    4  = 4273000000.0    # Inventory
    6  = 14418000000.0   # Total current assets
    7  = 23677000000.0   # Property, Plant, Equpment (Net)
    8  = 13376000000.0   # Long-Term Investments
    9  = 2261000000.0    # Goodwill and Intangible Assets (Total)
    10 = 1764000000.0    # Long-term assets (Other)
    11 = 41078000000.0   # Total non-current assets
    12 = 55496000000.0   # Total Assets
    13 = 8246000000.0    # Accounts Payable
    14 = 1149000000.0    # Accrued Expenses
    15 = 516000000.0     # Current Part of Debt
    18 = 1520000000.0    # Other current liabilities
    19 = 11431000000.0   # Total current liabilities
    20 = 14906000000.0   # Long Term Debt (Total)
    22 = 1351000000.0    # Pension and Post-Retirement Liabilities
    23 = 5547000000.0    # Long Term Tax Liability (Deferred)
    25 = 1804000000.0    # Non-current Liabilities (Other)
    26 = 23608000000.0   # Total non-current liabilities
    27 = 35039000000.0   # Total liabilities
    28 = 20420000000.0   # Additional Paid In Capital
    30 = 6000000.0       # Common Stock (Net)
    31 = 15449000000.0   # Retained Earnings
    33 = 17116000000.0   # Treasury Stock
    35 = 20457000000.0   # Common Equity (Total)
    36 = 20457000000.0   # Shareholders Equity (Total)
    37 = 55496000000.0   # Shareholders Equity and Liabilities (Total)
    38 = 437867000.0     # Shares (Common)
    39 = 18196000000.0   # Shareholders Equity (Tangible)
}

foreach ($row in $rowValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $rowValues[$row]
}
